$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used to restore the original (border-only) cell style
# after forcing a percent-looking string to stay literal text via a leading quote.
$xlPasteFormats = -4122

$ws.Range("E2").Value = "2026-02-09 23:18:26"
$ws.Range("I2").Value = "7.4 mm"
$ws.Range("E3").Value = "2026-02-09 23:18:29"
$ws.Range("H3").Value = "'97%"
$ws.Range("G3").Copy() | Out-Null
$ws.Range("H3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I3").Value = "6.3 mm"
$ws.Range("E4").Value = "2026-02-09 23:18:31"
$ws.Range("E5").Value = "2026-02-09 23:18:34"
$ws.Range("I5").Value = "3.4 mm"
$ws.Range("E6").Value = "2026-02-09 23:18:36"
$ws.Range("E7").Value = "2026-02-09 23:18:38"
$ws.Range("H7").Value = "'60%"
$ws.Range("G7").Copy() | Out-Null
$ws.Range("H7").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E8").Value = "2026-02-09 23:18:41"
$ws.Range("L8").Value = "41.4 km/h - 248º 22:55 TU"
$ws.Range("E9").Value = "2026-02-09 23:18:43"
$ws.Range("O9").Value = "8.4 °C"
$ws.Range("E10").Value = "2026-02-09 23:18:46"
$ws.Range("E11").Value = "2026-02-09 23:18:48"
$ws.Range("H11").Value = "'84%"
$ws.Range("G11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I11").Value = "1.3 mm"
$ws.Range("E12").Value = "2026-02-09 23:18:50"
$ws.Range("E13").Value = "2026-02-09 23:18:52"
$ws.Range("H13").Value = "'79%"
$ws.Range("G13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I13").Value = "1.9 mm"
$ws.Range("E14").Value = "2026-02-09 23:18:55"
$ws.Range("O14").Value = "10.6 °C"
$ws.Range("E15").Value = "2026-02-09 23:18:57"
$ws.Range("E16").Value = "2026-02-09 23:19:00"
$ws.Range("H16").Value = "'75%"
$ws.Range("G16").Copy() | Out-Null
$ws.Range("H16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I16").Value = "4.0 mm"
$ws.Range("E17").Value = "2026-02-09 23:19:02"
$ws.Range("H17").Value = "'83%"
$ws.Range("G17").Copy() | Out-Null
$ws.Range("H17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E18").Value = "2026-02-09 23:19:04"
$ws.Range("E19").Value = "2026-02-09 23:19:07"
$ws.Range("I19").Value = "1.3 mm"
$ws.Range("O19").Value = "4.3 °C"
$ws.Range("E20").Value = "2026-02-09 23:19:09"
$ws.Range("I20").Value = "1.8 mm"
$ws.Range("E21").Value = "2026-02-09 23:19:12"
$ws.Range("H21").Value = "'80%"
$ws.Range("G21").Copy() | Out-Null
$ws.Range("H21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I21").Value = "1.6 mm"
$ws.Range("E22").Value = "2026-02-09 23:19:14"
$ws.Range("O22").Value = "-4.8 °C"
$ws.Range("E23").Value = "2026-02-09 23:19:16"
$ws.Range("I23").Value = "5.9 mm"
$ws.Range("O23").Value = "-3.5 °C"
$ws.Range("E24").Value = "2026-02-09 23:19:19"
$ws.Range("I24").Value = "2.8 mm"
$ws.Range("E25").Value = "2026-02-09 23:19:21"
$ws.Range("I25").Value = "2.8 mm"
$ws.Range("E26").Value = "2026-02-09 23:19:24"
$ws.Range("I26").Value = "0.1 mm"
$ws.Range("E27").Value = "2026-02-09 23:19:26"
$ws.Range("I27").Value = "2.9 mm"
$ws.Range("E28").Value = "2026-02-09 23:19:29"
$ws.Range("H28").Value = "'81%"
$ws.Range("G28").Copy() | Out-Null
$ws.Range("H28").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I28").Value = "0.5 mm"
$ws.Range("E29").Value = "2026-02-09 23:19:31"
$ws.Range("E30").Value = "2026-02-09 23:19:33"
$ws.Range("E31").Value = "2026-02-09 23:19:36"
$ws.Range("E32").Value = "2026-02-09 23:19:38"
$ws.Range("I32").Value = "2.2 mm"
$ws.Range("E33").Value = "2026-02-09 23:19:41"
$ws.Range("H33").Value = "'79%"
$ws.Range("G33").Copy() | Out-Null
$ws.Range("H33").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I33").Value = "2.0 mm"
$ws.Range("E34").Value = "2026-02-09 23:19:43"
$ws.Range("H34").Value = "'77%"
$ws.Range("G34").Copy() | Out-Null
$ws.Range("H34").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I34").Value = "2.0 mm"
$ws.Range("O34").Value = "-1.0 °C"
$ws.Range("E35").Value = "2026-02-09 23:19:45"
$ws.Range("M35").Value = "10.1 °C 22:58 TU"
$ws.Range("O35").Value = "5.7 °C"
$ws.Range("E36").Value = "2026-02-09 23:19:48"
$ws.Range("H36").Value = "'80%"
$ws.Range("G36").Copy() | Out-Null
$ws.Range("H36").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E37").Value = "2026-02-09 23:19:50"
$ws.Range("I37").Value = "0.6 mm"
$ws.Range("J37").Value = "1007.9 hPa"
$ws.Range("E38").Value = "2026-02-09 23:19:52"
$ws.Range("I38").Value = "0.4 mm"
$ws.Range("E39").Value = "2026-02-09 23:19:55"
$ws.Range("H39").Value = "'78%"
$ws.Range("G39").Copy() | Out-Null
$ws.Range("H39").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I39").Value = "0.2 mm"
$ws.Range("E40").Value = "2026-02-09 23:19:57"
$ws.Range("I40").Value = "2.2 mm"
$ws.Range("O40").Value = "4.9 °C"
$ws.Range("E41").Value = "2026-02-09 23:19:59"
$ws.Range("H41").Value = "'60%"
$ws.Range("G41").Copy() | Out-Null
$ws.Range("H41").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J41").Value = "1007.6 hPa"
$ws.Range("O41").Value = "12.8 °C"
$ws.Range("E42").Value = "2026-02-09 23:20:02"
$ws.Range("E43").Value = "2026-02-09 23:20:04"
$ws.Range("I43").Value = "0.6 mm"
$ws.Range("E44").Value = "2026-02-09 23:20:07"
$ws.Range("I44").Value = "3.4 mm"
$ws.Range("E45").Value = "2026-02-09 23:20:09"
$ws.Range("G45").Value = "2 cm"
$ws.Range("H45").Value = "'85%"
$ws.Range("G45").Copy() | Out-Null
$ws.Range("H45").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("I45").Value = "3.2 mm"
$ws.Range("J45").Value = "1007.0 hPa"
$ws.Range("O45").Value = "4.2 °C"
$ws.Range("E46").Value = "2026-02-09 23:20:11"
